# Update the "Förändrad" date column (C) for rows 2-10 from 45174 to 45175
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 45175
